# Guidance - Instructions.docx edit script
# Applies the 4 changes described by the commit diff:
#  1. Remove the stray "_GoBack" bookmark that sits after "ataScienceScotland".
#  2. Reword "Use the Load new weights 2010 SAS project ..." sentence, splitting
#     it into two runs (quoting "Load new weights" with curly quotes).
#  3. Re-insert the "_GoBack" bookmark at the start of the "Run Prog. 3 ..."
#     paragraph (that's where Word's cursor was left when the doc was saved).
#  4. Reorder the "OLD instructions (replaced ... in December 2022)." sentence
#     into "OLD instructions (replaced in December 2022 ... above).", split into
#     three runs, but only in the first occurrence ("How it was developed...").
#
# Technique note: when two adjacent runs end up with identical formatting,
# a plain Range.Text/InsertAfter edit gets silently re-coalesced into a single
# run on save. Cutting a sub-range and pasting it back at the split point
# reliably creates (and keeps) a genuine run boundary while carrying the
# original run formatting along with it - so every multi-run split below is
# built that way, working from the right-most boundary back to the left so
# that later pastes never re-absorb an already created earlier boundary.

$d = $word.ActiveDocument

# --- 1. Remove the old "_GoBack" bookmark (was after "ataScienceScotland") ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- 2. Reword the "Load new weights" sentence, split into two runs ---
$q1 = [char]0x201C
$q2 = [char]0x201D

$lwOld = "Use the Load new weights 2010 SAS project (update the name) to import the weights file into SAS."
$lwPart1 = "Use the " + $q1 + "Load new weights" + $q2 + " SAS project"
$lwPart2 = " to import the weights file into SAS."

$lwRng = $d.Content
$lwFound = $lwRng.Find.Execute($lwOld, $true, $false, $false, $false, $false, $true, 1, $false, ($lwPart1 + $lwPart2), 2)

if ($lwFound) {
    # Cut the second half and paste it right back after the first half - this
    # forces a genuine run boundary between the two pieces.
    $lwCut = $d.Content
    $lwCut.Find.Execute($lwPart2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $lwCut.Cut()

    $lwPaste = $d.Content
    $lwPaste.Find.Execute($lwPart1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $lwPaste.Collapse(0)
    $lwPaste.Paste()
}

# --- 3. Re-insert the "_GoBack" bookmark at the start of "Run Prog. 3 ..." ---
$runPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Run Prog. 3*") {
        $runPara = $cand
        break
    }
}
if ($runPara -ne $null) {
    $runStart = $runPara.Range.Start
    $bmRange = $d.Range($runStart, $runStart)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# --- 4. Reorder the first "OLD instructions (...)" sentence, split into 3 runs ---
$oldSentence = "OLD instructions (replaced by the Running Instructions above in December 2022)."
$oiPart1 = "OLD instructions (replaced"
$oiPart2 = " in December 2022"
$oiPart3 = " by the Running Instructions above)."

function Find-OiParagraph($doc) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $cand = $doc.Paragraphs.Item($i)
        if ($cand.Range.Text -like "How it was developed*OLD instructions*") {
            return $cand
        }
    }
    return $null
}

$oiPara = Find-OiParagraph $d
if ($oiPara -ne $null) {
    $oiScopedRng = $oiPara.Range
    $oiFound = $oiScopedRng.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, ($oiPart1 + $oiPart2 + $oiPart3), 2)

    if ($oiFound) {
        # Right-most boundary first: cut part3, paste right after part2.
        $oiPara = Find-OiParagraph $d
        $oiCut3 = $oiPara.Range
        $oiCut3.Find.Execute($oiPart3, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $oiCut3.Cut()

        $oiPara = Find-OiParagraph $d
        $oiPaste3 = $oiPara.Range
        $oiPaste3.Find.Execute($oiPart2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $oiPaste3.Collapse(0)
        $oiPaste3.Paste()

        # Left boundary next: cut part2+part3, paste right after part1.
        $oiPara = Find-OiParagraph $d
        $oiCut23 = $oiPara.Range
        $oiCut23.Find.Execute(($oiPart2 + $oiPart3), $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $oiCut23.Cut()

        $oiPara = Find-OiParagraph $d
        $oiPaste23 = $oiPara.Range
        $oiPaste23.Find.Execute($oiPart1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
        $oiPaste23.Collapse(0)
        $oiPaste23.Paste()
    }
}

Write-Host "Edit script complete."
